# "PC1: changed 2  A.xlsx"
# - Sheet1!A2 shared string "pc1" -> "1-PC"
# - Sheet1!A3 shared string "pc2" -> "2-PC"
# - Active selection moved from A3 to F10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1-PC"
$ws.Range("A3").Value = "2-PC"

$ws.Range("F10").Select()
